# Auto-generated script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '56.978.21'
Set-TextValue $ws.Range('E2') '  +2.23%  '
Set-TextValue $ws.Range('D3') '3.010.46'
Set-TextValue $ws.Range('E3') '  +1.71%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.08%  '
Set-TextValue $ws.Range('D5') '514.99'
Set-TextValue $ws.Range('E5') '  +4.62%  '
Set-TextValue $ws.Range('D6') '139.67'
Set-TextValue $ws.Range('E6') '  +4.85%  '
Set-TextValue $ws.Range('D7') '0.999'
Set-TextValue $ws.Range('E7') '  +0.04%  '
Set-TextValue $ws.Range('E8') '  +3.11%  '
Set-TextValue $ws.Range('D9') '7.55'
Set-TextValue $ws.Range('E9') '  +5.98%  '
Set-TextValue $ws.Range('E10') '  +7.07%  '
Set-TextValue $ws.Range('E11') '  +3.27%  '
Set-TextValue $ws.Range('E12') '  +2.69%  '
Set-TextValue $ws.Range('D13') '3.524.21'
Set-TextValue $ws.Range('E13') '  +1.72%  '
Set-TextValue $ws.Range('D14') '25.75'
Set-TextValue $ws.Range('E14') '  +4.08%  '
Set-TextValue $ws.Range('D15') '0.0000158'
Set-TextValue $ws.Range('E15') '  +11.67%  '
Set-TextValue $ws.Range('D16') '56.983.81'
Set-TextValue $ws.Range('E16') '  +2.37%  '
Set-TextValue $ws.Range('D17') '3.002.17'
Set-TextValue $ws.Range('E17') '  +1.66%  '
Set-TextValue $ws.Range('D18') '5.95'
Set-TextValue $ws.Range('E18') '  +4.78%  '
Set-TextValue $ws.Range('D19') '12.59'
Set-TextValue $ws.Range('E19') '  +3.62%  '
Set-TextValue $ws.Range('D20') '7.90'
Set-TextValue $ws.Range('E20') '  +4.12%  '
Set-TextValue $ws.Range('D21') '328.63'
Set-TextValue $ws.Range('E21') '  +3.39%  '
Set-TextValue $ws.Range('E22') '  -0.02%  '
Set-TextValue $ws.Range('D23') '0.484'
Set-TextValue $ws.Range('E23') '  +4.89%  '
Set-TextValue $ws.Range('D24') '63.63'
Set-TextValue $ws.Range('E24') '  +5.56%  '
Set-TextValue $ws.Range('D25') '0.172'
Set-TextValue $ws.Range('E25') '  +6.84%  '
Set-TextValue $ws.Range('D26') '0.998'
Set-TextValue $ws.Range('E26') '  -0.28%  '
Set-TextValue $ws.Range('E27') '  +8.41%  '
Set-TextValue $ws.Range('D28') '6.66'
Set-TextValue $ws.Range('E28') '  +2.65%  '
Set-TextValue $ws.Range('D29') '7.14'
Set-TextValue $ws.Range('E29') '  +8.19%  '
Set-TextValue $ws.Range('E30') '  +7.24%  '
Set-TextValue $ws.Range('D31') '1.82'
Set-TextValue $ws.Range('E31') '  +6.67%  '
Set-TextValue $ws.Range('D32') '20.69'
Set-TextValue $ws.Range('E32') '  +5.90%  '
Set-TextValue $ws.Range('D33') '157.39'
Set-TextValue $ws.Range('E33') '  +5.40%  '
Set-TextValue $ws.Range('E34') '  +4.61%  '
Set-TextValue $ws.Range('D35') '5.74'
Set-TextValue $ws.Range('E35') '  +0.86%  '
Set-TextValue $ws.Range('E36') '  -1.69%  '
Set-TextValue $ws.Range('B37') 'Hedera'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D37') '0.0681'
Set-TextValue $ws.Range('E37') '  +3.72%  '
Set-TextValue $ws.Range('B38') 'EnergySwap'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D38') '24.29'
Set-TextValue $ws.Range('E38') '  +3.56%  '
Set-TextValue $ws.Range('D39') '3.042.08'
Set-TextValue $ws.Range('E39') '  +1.73%  '
Set-TextValue $ws.Range('D40') '37.23'
Set-TextValue $ws.Range('E40') '  +2.06%  '
Set-TextValue $ws.Range('D41') '1.00'
Set-TextValue $ws.Range('E41') '  +0.19%  '
Set-TextValue $ws.Range('D42') '2.291.58'
Set-TextValue $ws.Range('E42') '  +8.34%  '
Set-TextValue $ws.Range('D43') '0.649'
Set-TextValue $ws.Range('E43') '  +2.86%  '
Set-TextValue $ws.Range('B44') 'Filecoin'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D44') '3.70'
Set-TextValue $ws.Range('E44') '  +4.55%  '
Set-TextValue $ws.Range('B45') 'Stacks'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D45') '1.43'
Set-TextValue $ws.Range('E45') '  +2.86%  '
Set-TextValue $ws.Range('E46') '  +1.62%  '
Set-TextValue $ws.Range('E47') '  +8.03%  '
Set-TextValue $ws.Range('E48') '  +2.93%  '
Set-TextValue $ws.Range('E49') '  +6.30%  '
Set-TextValue $ws.Range('D50') '19.31'
Set-TextValue $ws.Range('E50') '  -0.05%  '
Set-TextValue $ws.Range('D51') '0.0879'
